$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (shared strings with rich-text runs) ---
# A8: "Volume 30   Number  46" -> "...  47"
$volCell = $ws.Range("A8")
$volChars = $volCell.Characters(21, 2)
$volChars.Text = "47"

# C9: "Report Covering the Week  11/13/2023  Through  11/19/2023"
#  -> "...  11/20/2023  Through  11/26/2023"
$weekCell = $ws.Range("C9")
$weekChars1 = $weekCell.Characters(27, 10)
$weekChars1.Text = "11/20/2023"
$weekChars2 = $weekCell.Characters(48, 10)
$weekChars2.Text = "11/26/2023"

# --- Update Crime Complaints table (rows 14-30) ---
$ws.Range("C14").Copy($ws.Range("F14"))
$ws.Range("H14").Value = -100
$ws.Range("M14").Value = -58.823529411764
$ws.Range("G14").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("G14").Copy($ws.Range("F15"))
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = -9.523809523809
$ws.Range("L15").Value = -17.391304347826
$ws.Range("M15").Value = -5
$ws.Range("N15").Value = -71.212121212121
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 11.764705882352
$ws.Range("I16").Value = 201
$ws.Range("J16").Value = 196
$ws.Range("K16").Value = 2.551020408163
$ws.Range("L16").Value = 58.267716535433
$ws.Range("M16").Value = -16.942148760330
$ws.Range("N16").Value = -89.802130898021
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -50
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -17.142857142857
$ws.Range("I17").Value = 378
$ws.Range("J17").Value = 338
$ws.Range("K17").Value = 11.834319526627
$ws.Range("L17").Value = 33.098591549295
$ws.Range("M17").Value = 35.483870967741
$ws.Range("N17").Value = -65.791855203619
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 169
$ws.Range("K18").Value = -0.591715976331
$ws.Range("L18").Value = 25.373134328358
$ws.Range("M18").Value = -41.463414634146
$ws.Range("N18").Value = -88.373702422145
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -19.444444444444
$ws.Range("I19").Value = 362
$ws.Range("J19").Value = 397
$ws.Range("K19").Value = -8.816120906801
$ws.Range("L19").Value = 41.40625
$ws.Range("M19").Value = 31.636363636363
$ws.Range("N19").Value = -32.209737827715
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 153
$ws.Range("J20").Value = 122
$ws.Range("K20").Value = 25.409836065573
$ws.Range("L20").Value = 62.765957446808
$ws.Range("M20").Value = 29.661016949152
$ws.Range("N20").Value = -78.017241379310
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 1288
$ws.Range("J21").Value = 1254
$ws.Range("K21").Value = 2.711323763955
$ws.Range("L21").Value = 38.643702906350
$ws.Range("M21").Value = 4.038772213247
$ws.Range("N21").Value = -78.136139874384
$ws.Range("G14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 100
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -16.129032258064
$ws.Range("L22").Value = 52.941176470588
$ws.Range("C23").Value = 1
$ws.Range("G14").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 4
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 107
$ws.Range("J23").Value = 97
$ws.Range("K23").Value = 10.309278350515
$ws.Range("L23").Value = 27.380952380952
$ws.Range("M23").Value = 64.615384615384
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -26.315789473684
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 702
$ws.Range("J24").Value = 804
$ws.Range("K24").Value = -12.686567164179
$ws.Range("L24").Value = 12.680577849117
$ws.Range("M24").Value = -8.355091383812
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 18.421052631578
$ws.Range("I25").Value = 513
$ws.Range("J25").Value = 445
$ws.Range("K25").Value = 15.280898876404
$ws.Range("L25").Value = 51.327433628318
$ws.Range("M25").Value = -30.204081632653
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 33
$ws.Range("K26").Value = -5.714285714285
$ws.Range("L26").Value = -2.941176470588
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -9.722222222222
$ws.Range("G14").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -85.714285714285
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = -39.473684210526
$ws.Range("M28").Value = -65.671641791044
$ws.Range("N28").Value = -90.456431535269
$ws.Range("G14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -80
$ws.Range("J29").Value = 33
$ws.Range("K29").Value = -36.363636363636
$ws.Range("M29").Value = -60.377358490566
$ws.Range("N29").Value = -90.322580645161
$ws.Range("G14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = -33.333333333333
$ws.Range("L30").Value = -20
